$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell C1
$ws.Range("C1").Value = "Haute PAGE 32 OF  33"

# Append three more rows identical to row 14 (Potential Usage of Vulnerable Log4J)
$ws.Range("A15").Value = "Potential Usage of Vulnerable Log4J"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "1"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "Informations"

$ws.Range("A16").Value = "Potential Usage of Vulnerable Log4J"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "1"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "Informations"

$ws.Range("A17").Value = "Potential Usage of Vulnerable Log4J"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "1"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "Informations"
